$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 322
$ws.Range("I9").Value = 143.33333
$ws.Range("J9").Value = 500.66666
$ws.Range("K9").Value = 143.33333
$ws.Range("L9").Value = 500.66666
$ws.Range("M9").Value = 25.66667000000001
$ws.Range("N9").Value = -838.66666
$ws.Range("H33").Value = 256.8421
$ws.Range("I33").Value = 271.7647
$ws.Range("K33").Value = 271.7647
$ws.Range("M33").Value = -42.7647
$ws.Range("H62").Value = 4538.05
$ws.Range("I62").Value = 3967
$ws.Range("J62").Value = 5109.1
$ws.Range("K62").Value = 3967
$ws.Range("L62").Value = 5109.1
$ws.Range("M62").Value = -3343
$ws.Range("N62").Value = -6357.1
$ws.Range("H65").Value = 4538.05
$ws.Range("I65").Value = 3967
$ws.Range("J65").Value = 5109.1
$ws.Range("K65").Value = 19835
$ws.Range("L65").Value = 25545.5
$ws.Range("M65").Value = -16715
$ws.Range("N65").Value = -31785.5
$ws.Range("H80").Value = 11075934
$ws.Range("I80").Value = 444.66666
$ws.Range("J80").Value = 18743580
$ws.Range("K80").Value = 1333.99998
$ws.Range("L80").Value = 56230740
$ws.Range("M80").Value = -335.9999800000001
$ws.Range("N80").Value = -56232736
$ws.Range("H83").Value = 11075934
$ws.Range("I83").Value = 444.66666
$ws.Range("J83").Value = 18743580
$ws.Range("K83").Value = 4001.99994
$ws.Range("L83").Value = 168692220
$ws.Range("M83").Value = 990.0000600000003
$ws.Range("N83").Value = -168702204
$ws.Range("H96").Value = 125000760
$ws.Range("I96").Value = 125000760
$ws.Range("K96").Value = 375002280
$ws.Range("M96").Value = -375000907
$ws.Range("H98").Value = 857.7778
$ws.Range("I98").Value = 840
$ws.Range("K98").Value = 840
$ws.Range("M98").Value = 658
$ws.Range("H113").Value = 90913590
$ws.Range("I113").Value = 200001400
$ws.Range("J113").Value = 7081
$ws.Range("K113").Value = 200001400
$ws.Range("L113").Value = 7081
$ws.Range("M113").Value = -199998146
$ws.Range("N113").Value = -13589
$ws.Range("H122").Value = 857.7778
$ws.Range("I122").Value = 840
$ws.Range("K122").Value = 2520
$ws.Range("M122").Value = -70
$ws.Range("H129").Value = 1638.1936
$ws.Range("J129").Value = 1792.3928
$ws.Range("L129").Value = 5377.178400000001
$ws.Range("N129").Value = -15377.1784
$ws.Range("H137").Value = 130932.484
$ws.Range("I137").Value = 144410.97
$ws.Range("K137").Value = 433232.91
$ws.Range("M137").Value = -430682.91
$ws.Range("H138").Value = 4303.273
$ws.Range("I138").Value = 4660
$ws.Range("J138").Value = 4239.5713
$ws.Range("K138").Value = 13980
$ws.Range("L138").Value = 12718.7139
$ws.Range("M138").Value = -8840
$ws.Range("N138").Value = -22998.7139
$ws.Range("H141").Value = 1565.55
$ws.Range("I141").Value = 1261.9706
$ws.Range("J141").Value = 3285.8333
$ws.Range("K141").Value = 3785.9118
$ws.Range("L141").Value = 9857.499899999999
$ws.Range("M141").Value = 1394.0882
$ws.Range("N141").Value = -20217.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 200
$ws.Range("I25").Value = 200
$ws.Range("K25").Value = 200
$ws.Range("M25").Value = 202
$ws.Range("H32").Value = 9938.608
$ws.Range("I32").Value = 7180.793
$ws.Range("K32").Value = 7180.793
$ws.Range("M32").Value = -6893.793
$ws.Range("H35").Value = 1616.5
$ws.Range("I35").Value = 1616.5
$ws.Range("K35").Value = 1616.5
$ws.Range("M35").Value = -1210.5
$ws.Range("H45").Value = 2739.516
$ws.Range("I45").Value = 2495.0952
$ws.Range("J45").Value = 3252.8
$ws.Range("K45").Value = 2495.0952
$ws.Range("L45").Value = 3252.8
$ws.Range("M45").Value = -2118.0952
$ws.Range("N45").Value = -4006.8
$ws.Range("H97").Value = 439.0909
$ws.Range("I97").Value = 278.625
$ws.Range("K97").Value = 278.625
$ws.Range("M97").Value = 217.375
$ws.Range("H137").Value = 47618
$ws.Range("I137").Value = 43000
$ws.Range("J137").Value = 48772.5
$ws.Range("K137").Value = 43000
$ws.Range("L137").Value = 48772.5
$ws.Range("M137").Value = -37900
$ws.Range("N137").Value = -58972.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3819.6155
$ws.Range("I134").Value = 3867.5
$ws.Range("K134").Value = 11602.5
$ws.Range("M134").Value = -9067.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 40.857143
$ws.Range("I7").Value = 40
$ws.Range("J7").Value = 41.2
$ws.Range("K7").Value = 40
$ws.Range("L7").Value = 41.2
$ws.Range("M7").Value = 73
$ws.Range("N7").Value = -267.2
$ws.Range("H58").Value = 20120.889
$ws.Range("I58").Value = 1672.6923
$ws.Range("J58").Value = 37251.355
$ws.Range("K58").Value = 1672.6923
$ws.Range("L58").Value = 37251.355
$ws.Range("M58").Value = -1469.6923
$ws.Range("N58").Value = -37657.355
$ws.Range("H97").Value = 31464
$ws.Range("J97").Value = 31464
$ws.Range("L97").Value = 31464
$ws.Range("N97").Value = -33446
$ws.Range("H132").Value = 41669636
$ws.Range("I132").Value = 66668924
$ws.Range("J132").Value = 4160.778
$ws.Range("K132").Value = 200006772
$ws.Range("L132").Value = 12482.334
$ws.Range("M132").Value = -200004242
$ws.Range("N132").Value = -17542.334
$ws.Range("H136").Value = 20120.889
$ws.Range("I136").Value = 1672.6923
$ws.Range("J136").Value = 37251.355
$ws.Range("K136").Value = 5018.0769
$ws.Range("L136").Value = 111754.065
$ws.Range("M136").Value = -2468.0769
$ws.Range("N136").Value = -116854.065

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1367.3823
$ws.Range("I5").Value = 919.26086
$ws.Range("K5").Value = 2757.78258
$ws.Range("M5").Value = -2645.78258
$ws.Range("H131").Value = 800.86
$ws.Range("I131").Value = 476.5
$ws.Range("J131").Value = 836.9
$ws.Range("K131").Value = 1429.5
$ws.Range("L131").Value = 2510.7
$ws.Range("M131").Value = 3610.5
$ws.Range("N131").Value = -12590.7
$ws.Range("H135").Value = 1367.3823
$ws.Range("I135").Value = 919.26086
$ws.Range("K135").Value = 8273.347739999999
$ws.Range("M135").Value = -5738.347739999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6535
$ws.Range("I122").Value = 7686.25
$ws.Range("K122").Value = 23058.75
$ws.Range("M122").Value = -20608.75
$ws.Range("H126").Value = 3476.0688
$ws.Range("I126").Value = 2587.4666
$ws.Range("J126").Value = 4428.143
$ws.Range("K126").Value = 7762.399800000001
$ws.Range("L126").Value = 13284.429
$ws.Range("M126").Value = -5292.399800000001
$ws.Range("N126").Value = -18224.429
$ws.Range("H132").Value = 4117086.8
$ws.Range("I132").Value = 6688778
$ws.Range("K132").Value = 20066334
$ws.Range("M132").Value = -20063804

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2930
$ws.Range("I7").Value = 2550
$ws.Range("K7").Value = 2550
$ws.Range("M7").Value = -2438
$ws.Range("H40").Value = 4442.6523
$ws.Range("I40").Value = 4235.5
$ws.Range("K40").Value = 4235.5
$ws.Range("M40").Value = -4099.5
$ws.Range("H61").Value = 4017
$ws.Range("I61").Value = 2231.1667
$ws.Range("J61").Value = 8303
$ws.Range("K61").Value = 2231.1667
$ws.Range("L61").Value = 8303
$ws.Range("M61").Value = -2029.1667
$ws.Range("N61").Value = -8707
$ws.Range("H113").Value = 4017
$ws.Range("I113").Value = 2231.1667
$ws.Range("J113").Value = 8303
$ws.Range("K113").Value = 2231.1667
$ws.Range("L113").Value = 8303
$ws.Range("M113").Value = -61.16670000000022
$ws.Range("N113").Value = -12643
$ws.Range("H126").Value = 2930
$ws.Range("I126").Value = 2550
$ws.Range("K126").Value = 7650
$ws.Range("M126").Value = -5180
$ws.Range("H132").Value = 243635.8
$ws.Range("I132").Value = 356667.5
$ws.Range("J132").Value = 3443.4375
$ws.Range("K132").Value = 1070002.5
$ws.Range("L132").Value = 10330.3125
$ws.Range("M132").Value = -1067472.5
$ws.Range("N132").Value = -15390.3125
$ws.Range("H136").Value = 2213.087
$ws.Range("I136").Value = 2222.7727
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 6668.3181
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -4118.3181
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3125
$ws.Range("J15").Value = 3125
$ws.Range("L15").Value = 3125
$ws.Range("N15").Value = -3701
$ws.Range("H20").Value = 4747.8335
$ws.Range("J20").Value = 5097.4
$ws.Range("L20").Value = 5097.4
$ws.Range("N20").Value = -5577.4
$ws.Range("H113").Value = 796.8929000000001
$ws.Range("I113").Value = 1035.6875
$ws.Range("J113").Value = 478.5
$ws.Range("K113").Value = 3107.0625
$ws.Range("L113").Value = 1435.5
$ws.Range("M113").Value = -937.0625
$ws.Range("N113").Value = -5775.5
$ws.Range("H136").Value = 38235580
$ws.Range("I136").Value = 54331296
$ws.Range("K136").Value = 162993888
$ws.Range("M136").Value = -162991338
